$d = $word.ActiveDocument

# Locate the "Write Up" title paragraph so we insert relative to it rather than
# assuming a fixed paragraph index.
$titleRange = $d.Content
[void]$titleRange.Find.Execute("Write Up", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$titlePara = $titleRange.Paragraphs.Item(1)

# Insert three new paragraphs right after the title:
#   1) intro paragraph (Normal style)
#   2) lead-in paragraph (Normal style)
#   3) "13 Door System" heading (Heading1 style)
$titlePara.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($titlePara.Index + 1)
$p2.Style = "Normal"
$p2.Range.Text = "In this tutorial, we will be starting to build our door system. In the Dungeon Crawler game, we will want to have a variety of different doors, which can be coded to take the hero to various areas, in the game."

$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item($titlePara.Index + 2)
$p3.Style = "Normal"
$p3.Range.Text = "So, if you are ever in need of a door system to be generated in your own game, and would like to know how to go about it, then please join us for our brand-new article this week entitled:"

$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item($titlePara.Index + 3)
$p4.Style = "Heading1"
$p4.Range.Text = "13 Door System"
